# Updated cryptos list values (prices + 1h volume %) and one row swap (Toncoin <-> InjectiveProtocol)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.586.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.283.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "94.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.623.04"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.844"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.283.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.514.82"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.17"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -12.34%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.47"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.39"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.83"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.23%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.07%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.235"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +15.69%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.50"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.21"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.51"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.501.58"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.06%  "
